# Extended Config to 2K, refactors Tasks and Config
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "System State machine" sheet: add state-machine state names in A2:A11
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("System State machine")

$ws3.Range("A2").Value  = "Initialisation"
$ws3.Range("A3").Value  = "Set Date/Time"
$ws3.Range("A4").Value  = "Enter Task  Config"
$ws3.Range("A5").Value  = "Enter Face Config"
$ws3.Range("A6").Value  = "IDLE_STATE"
$ws3.Range("A7").Value  = "MOVE_STATE"
$ws3.Range("A8").Value  = "NEW_FACE_DETECTED"
$ws3.Range("A9").Value  = "SLEEP_STATE"
$ws3.Range("A10").Value = "CHANGE_TASK"
$ws3.Range("A11").Value = "END_TASK"

$ws3.Columns.Item(1).ColumnWidth = 26
$ws3.Range("A12").Select()

# ---------------------------------------------------------------------
# 2) New "Sheet1" worksheet, placed after "System State machine", with
#    the reworked / extended Config sizing tables.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws4.Name = "Sheet1"

$ws4.Range("A2").Value = "static uint16_t SYS_CONFIG_SIZE ;"
$ws4.Range("B2").Value = 896
$ws4.Range("C2").Formula = "=B2/4"
$ws4.Range("D2").Formula = "=B2*8"
$ws4.Range("E2").Formula = "=D2/32"

$ws4.Range("A3").Value = "static uint16_t SYS_CONFIG_ALL_SIZE;"
$ws4.Range("B3").Value = 904
$ws4.Range("A4").Value = "static uint16_t SYS_CONFIG_SIZE_TASK;"
$ws4.Range("B4").Value = 800
$ws4.Range("A5").Value = "static uint16_t SYS_CONFIG_SIZE_DODEC;"
$ws4.Range("B5").Value = 96

$ws4.Range("C3:C5").Formula = "=B3/4"
$ws4.Range("D3:D5").Formula = "=B3*8"
$ws4.Range("E3:E5").Formula = "=D3/32"

$ws4.Range("A8").Value = "Where is "
$ws4.Range("B8").Formula = "=B3-B2"

$ws4.Range("A10").Value = "static uint16_t SYS_CONFIG_SIZE ;"
$ws4.Range("B10").Value = 896
$ws4.Range("C10").Formula = "=B10/4"
$ws4.Range("D10").Formula = "=B10*8"
$ws4.Range("E10").Formula = "=D10/32"

$ws4.Range("A11").Value = "static uint16_t SYS_CONFIG_ALL_SIZE;"
$ws4.Range("B11").Value = 900
$ws4.Range("A12").Value = "static uint16_t SYS_CONFIG_SIZE_TASK;"
$ws4.Range("B12").Value = 800
$ws4.Range("A13").Value = "static uint16_t SYS_CONFIG_SIZE_DODEC;"
$ws4.Range("B13").Value = 96

$ws4.Range("C11:C13").Formula = "=B11/4"
$ws4.Range("D11:D13").Formula = "=B11*8"
$ws4.Range("E11:E13").Formula = "=D11/32"

$ws4.Range("A15").Value = "plud Checksum"
$ws4.Range("B15").Value = 4

$ws4.Columns.Item(1).ColumnWidth = 38
$ws4.Range("B16").Select()
$ws4.Activate()

# ---------------------------------------------------------------------
# 3) Config sheet: shared formula on D2 now only spans D2:D5
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Config")
$ws2.Range("D2:D5").Formula = "=C2*8"
